$d = $word.ActiveDocument

# --- helper pattern -------------------------------------------------
# This document's paragraphs all inherit "jc=center / Times New Roman"
# formatting automatically from the preceding paragraph when a new
# paragraph mark is inserted at the very end of the story, so we only
# need to set Bold / BoldBi / Superscript explicitly where the diff
# calls for them.
#
# Note: `$d.Content.End` always points just *before* the document's
# very last paragraph mark (standard Word behaviour), so re-reading
# `$d.Range($d.Content.End, $d.Content.End)` after every mutation gives
# the correct "end of story" insertion point.

# 1) blank separator paragraph after "Research for unit testing."
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()

# 2) "November 8th: " heading (bold)
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Font.Bold = $true
$p.Range.Font.BoldBi = $true

$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("November 8")
$r.Font.Bold = $true
$r.Font.BoldBi = $true

$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("th")
$r.Font.Bold = $true
$r.Font.BoldBi = $true
$r.Font.Superscript = $true

$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter(": ")
$r.Font.Bold = $true
$r.Font.BoldBi = $true

# 3) blank separator
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()

# 4) "Worked on User Interface."
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("Worked on User Interface.")

# 5) blank separator
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()

# 6) "November 15th:" heading (bold)
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.Font.Bold = $true
$p.Range.Font.BoldBi = $true

$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("November 15")
$r.Font.Bold = $true
$r.Font.BoldBi = $true

$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("th")
$r.Font.Bold = $true
$r.Font.BoldBi = $true
$r.Font.Superscript = $true

$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter(":")
$r.Font.Bold = $true
$r.Font.BoldBi = $true

# 7) blank separator
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()

# 8) "Worked on Friends List on dashboard."
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("Worked on Friends List on dashboard.")

# 9) blank separator
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()

# 10) "Updated Architecture Design doc."
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertParagraphAfter()
$r = $d.Range($d.Content.End, $d.Content.End)
$r.InsertAfter("Updated Architecture Design doc.")

Write-Output "applied"
